# Updated cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.143.05"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "'1.996.69"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'330.99"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("D7").Value = "'0.4983"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").Value = "'0.4203"
$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("D9").Value = "'54.73"
$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("D10").Value = "'0.08972"
$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("D11").Value = "'1.098"
$ws.Range("E11").Value = "  -2.26%  "

$ws.Range("D12").Value = "'23.22"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'1.996.93"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").Value = "'8.037"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").Value = "'6.444"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").Value = "'1.012"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'92.66"
$ws.Range("E17").Value = "  -3.37%  "

$ws.Range("D18").Value = "'0.00001109"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'0.06767"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").Value = "'19.62"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "'5.985"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").Value = "'29.164.05"
$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").Value = "'12.01"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").Value = "'2.233.85"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").Value = "'20.86"
$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").Value = "'157.19"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("D29").Value = "'6.321"
$ws.Range("E29").Value = "  -3.63%  "

$ws.Range("D30").Value = "'2.264"
$ws.Range("E30").Value = "  -2.81%  "

$ws.Range("D31").Value = "'127.51"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").Value = "'1.050"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").Value = "'0.09873"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").Value = "'1.533"
$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("D35").Value = "'5.832"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").Value = "'3.742"

$ws.Range("D37").Value = "'0.02427"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").Value = "'1.321"
$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("D39").Value = "'9.074"
$ws.Range("E39").Value = "  -5.50%  "

$ws.Range("D40").Value = "'0.06400"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").Value = "'0.6501"
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("D42").Value = "'11.51"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("D43").Value = "'0.1991"
$ws.Range("E43").Value = "  -3.71%  "

$ws.Range("D44").Value = "'1.011"
$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("D45").Value = "'0.6226"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("D46").Value = "'1.356"
$ws.Range("E46").Value = "  +6.53%  "

$ws.Range("D47").Value = "'13.37"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").Value = "'2.191"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("D50").Value = "'0.00000000335"
$ws.Range("E50").Value = "  +4.20%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'2.147"
$ws.Range("E51").Value = "  +10.44%  "
